$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (G2: valor faturado, H2: valor da comissão)
$ws.Range("G2").Value = 22468.0
$ws.Range("H2").Value = 2246.8

# Add new row 4 with order data
# Numeric-looking identifiers are entered as text (leading apostrophe keeps
# them as text without reformatting the cell's number format).
$ws.Range("A4").Value = "'4"
$ws.Range("B4").Value = "'123213"
$ws.Range("C4").Value = "Cliente Pedro 1"
$ws.Range("D4").Value = "'23"
$ws.Range("E4").Value = "TOTALMENTEFATURADO"
$ws.Range("F4").Value = 200001.0
$ws.Range("G4").Value = 200001.0
$ws.Range("H4").Value = 9000.045
